# Refresh the scraped cryptocurrency Price / Volume(1h) figures to the
# Sat May 11 10:36:07 UTC 2024 snapshot (GitHub Actions cron run).
# Also re-orders the Kaspa / Stacks rows (40 <-> 41) to match the
# refreshed coinranking.com market-cap ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value.
$updates = [ordered]@{
    "D2" = "60.777.61"
    "E2" = "  -3.61%  "
    "D3" = "2.903.41"
    "E3" = "  -4.30%  "
    "E4" = "  +0.05%  "
    "D5" = "586.71"
    "E5" = "  -1.18%  "
    "D6" = "144.23"
    "E6" = "  -6.27%  "
    "E7" = "  +0.06%  "
    "D8" = "0.502"
    "E8" = "  -2.51%  "
    "D9" = "2.902.33"
    "E9" = "  -4.22%  "
    "E10" = "  -3.81%  "
    "D11" = "0.143"
    "E11" = "  -5.57%  "
    "D12" = "0.443"
    "E12" = "  -4.40%  "
    "D13" = "0.0000225"
    "E13" = "  -3.81%  "
    "D14" = "33.42"
    "E14" = "  -6.35%  "
    "E15" = "  +1.57%  "
    "D16" = "3.385.20"
    "E16" = "  -4.29%  "
    "D17" = "60.774.07"
    "E17" = "  -3.57%  "
    "D18" = "6.66"
    "E18" = "  -6.13%  "
    "D19" = "2.906.65"
    "E19" = "  -4.16%  "
    "D20" = "427.49"
    "E20" = "  -5.84%  "
    "D21" = "13.50"
    "E21" = "  -5.64%  "
    "D22" = "0.680"
    "E22" = "  -2.58%  "
    "D23" = "7.08"
    "E23" = "  -5.85%  "
    "D24" = "80.78"
    "E24" = "  -2.87%  "
    "D25" = "10.91"
    "E25" = "  -4.48%  "
    "D26" = "2.22"
    "E26" = "  -5.07%  "
    "D27" = "11.88"
    "E27" = "  -4.50%  "
    "E28" = "  -0.06%  "
    "E29" = "  +0.05%  "
    "D30" = "2.19"
    "E30" = "  -4.00%  "
    "D31" = "7.21"
    "E31" = "  -4.04%  "
    "D32" = "2.60"
    "E32" = "  -3.51%  "
    "D33" = "26.33"
    "E33" = "  -4.73%  "
    "D34" = "0.107"
    "D35" = "0.0₃0859"
    "E35" = "  -0.48%  "
    "E36" = "  -3.56%  "
    "D37" = "5.56"
    "E37" = "  -6.04%  "
    "D38" = "3.03"
    "E38" = "  -3.92%  "
    "D39" = "49.47"
    "E39" = "  -1.95%  "
    "B40" = "Stacks"
    "C40" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "D40" = "1.99"
    "E40" = "  -5.88%  "
    "B41" = "Kaspa"
    "C41" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D41" = "0.123"
    "E41" = "  -6.06%  "
    "D42" = "8.57"
    "E42" = "  -5.88%  "
    "D43" = "0.295"
    "E43" = "  -5.12%  "
    "D44" = "41.27"
    "E44" = "  -7.44%  "
    "D45" = "0.0350"
    "E45" = "  -3.03%  "
    "D46" = "375.02"
    "E46" = "  -5.18%  "
    "D47" = "2.691.86"
    "E47" = "  -1.11%  "
    "D48" = "132.18"
    "E48" = "  -0.99%  "
    "D50" = "24.03"
    "E50" = "  -6.52%  "
    "E51" = "  -3.00%  "
}

# Values that look like plain numbers (e.g. "586.71") would otherwise be
# auto-converted to the Number type by the Value setter. The sheet stores
# these as plain text, so we briefly force a text NumberFormat on write and
# then restore the default (unstyled) cell style, exactly as typing the
# same text into Excel with the column pre-formatted as Text would.
$numericPattern = '^[+-]?[0-9]+(\.[0-9]+)?$'

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    if ($value -match $numericPattern) {
        $cell.NumberFormat = "@"
        $cell.Value = $value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $value
    }
}
